# Logboek.xlsx edit script
# Adds a "Dag" (Day) column to both the "Karsten" and "Douwe" worksheets,
# populating a few rows with the day of week, and switches the active tab
# from "Karsten" to "Douwe".

$wb = $excel.ActiveWorkbook

# --- Sheet "Karsten" ---
$wsK = $wb.Worksheets.Item("Karsten")
$wsK.Range("E1").Value = "Dag"
$wsK.Range("E2").Value = "Maandag"
$wsK.Range("E5").Value = "dinsdag"
$wsK.Range("E7").Value = "woensdag"

# Match header style (bold, centered, bordered) used by the rest of row 1
$wsK.Range("D1").Copy()
$wsK.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsK.Columns.Item(5).ColumnWidth = 9.1666666666666661

$wsK.Range("C4").Select()

# --- Sheet "Douwe" ---
$wsD = $wb.Worksheets.Item("Douwe")
$wsD.Range("E1").Value = "Dag"
$wsD.Range("E2").Value = "Maandag"
$wsD.Range("E4").Value = "dinsdag"
$wsD.Range("E6").Value = "woensdag"

$wsD.Range("D1").Copy()
$wsD.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsD.Range("D8").Select()

# Make "Douwe" the active sheet/tab
$wsD.Activate()
